$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prime style templates for the new rows by copying formats from row 120
$ws.Range("A120").Copy() | Out-Null
$ws.Range("A121:A122").PasteSpecial(-4122) | Out-Null
$ws.Range("E120").Copy() | Out-Null
$ws.Range("E121:E122").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- row 121 ----
$ws.Range("A121").Value = 120
$ws.Range("B121").Value = 'poland'
$ws.Range("C121").Value = 'ekstraklasa'
$ws.Range("D121").Value = '2023-2024'
$ws.Range("E121").Value = 45235.52083333334
$ws.Range("F121").Value = 'Slask Wroclaw'
$ws.Range("G121").Value = 2
$ws.Range("H121").Value = 'LKS Lodz'
$ws.Range("I121").Value = 1
$ws.Range("J121").Value = 1.76
$ws.Range("K121").Value = '29/10/2023 12:42'
$ws.Range("L121").Value = 1.73
$ws.Range("M121").Value = '05/11/2023 12:29'
$ws.Range("N121").Value = 3.75
$ws.Range("O121").Value = '29/10/2023 12:42'
$ws.Range("P121").Value = 3.68
$ws.Range("Q121").Value = '05/11/2023 12:29'
$ws.Range("R121").Value = 4.78
$ws.Range("S121").Value = '29/10/2023 12:42'
$ws.Range("T121").Value = 5.36
$ws.Range("U121").Value = '05/11/2023 12:29'
$ws.Range("V121").Value = 'https://www.betexplorer.com/football/poland/ekstraklasa/slask-wroclaw-lks-lodz/fBwrvKVc/'

# ---- row 122 ----
$ws.Range("A122").Value = 121
$ws.Range("B122").Value = 'poland'
$ws.Range("C122").Value = 'ekstraklasa'
$ws.Range("D122").Value = '2023-2024'
$ws.Range("E122").Value = 45235.625
$ws.Range("F122").Value = 'Rakow'
$ws.Range("G122").Value = 5
$ws.Range("H122").Value = 'Zaglebie'
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1.59
$ws.Range("K122").Value = '29/10/2023 17:43'
$ws.Range("L122").Value = 1.56
$ws.Range("M122").Value = '05/11/2023 14:34'
$ws.Range("N122").Value = 4.12
$ws.Range("O122").Value = '29/10/2023 17:43'
$ws.Range("P122").Value = 4.31
$ws.Range("Q122").Value = '05/11/2023 14:54'
$ws.Range("R122").Value = 5.08
$ws.Range("S122").Value = '29/10/2023 17:43'
$ws.Range("T122").Value = 5.88
$ws.Range("U122").Value = '05/11/2023 14:54'
$ws.Range("V122").Value = 'https://www.betexplorer.com/football/poland/ekstraklasa/rakow-czestochowa-zaglebie/bqpZut1o/'
